$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.597.29'
$ws.Range("E2").Value = '  +1.09%  '
$ws.Range("D3").Value = '1.562.04'
$ws.Range("E3").Value = '  +0.59%  '
$ws.Range("E4").Value = '  -0.40%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.27'
$ws.Range("E5").Value = '  +0.34%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.514'
$ws.Range("E6").Value = '  +5.79%  '
$ws.Range("E7").Value = '  -0.49%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '25.02'
$ws.Range("E8").Value = '  +6.84%  '
$ws.Range("E9").Value = '  +1.70%  '
$ws.Range("E10").Value = '  +0.72%  '
$ws.Range("E11").Value = '  +1.00%  '
$ws.Range("D12").Value = '1.784.11'
$ws.Range("E12").Value = '  +0.52%  '
$ws.Range("D13").Value = '1.565.63'
$ws.Range("E13").Value = '  +1.03%  '
$ws.Range("D14").Value = '28.619.14'
$ws.Range("E14").Value = '  +1.16%  '
$ws.Range("E15").Value = '  +1.53%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.65'
$ws.Range("E16").Value = '  +0.79%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.58'
$ws.Range("E17").Value = '  +2.02%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '227.77'
$ws.Range("E18").Value = '  +0.21%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.29'
$ws.Range("E19").Value = '  -0.10%  '
$ws.Range("E20").Value = '  +2.06%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.998'
$ws.Range("E21").Value = '  -0.44%  '
$ws.Range("E22").Value = '  +0.76%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.09'
$ws.Range("E23").Value = '  +3.43%  '
$ws.Range("E24").Value = '  +3.15%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.47'
$ws.Range("E25").Value = '  +2.42%  '
$ws.Range("E26").Value = '  +3.89%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '14.81'
$ws.Range("E27").Value = '  +0.41%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.27'
$ws.Range("E28").Value = '  +0.64%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.998'
$ws.Range("E29").Value = '  -0.46%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0459'
$ws.Range("E30").Value = '  -2.05%  '
$ws.Range("E31").Value = '  -0.99%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.19'
$ws.Range("E32").Value = '  +0.80%  '
$ws.Range("D33").Value = '1.405.67'
$ws.Range("E33").Value = '  +1.84%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.01'
$ws.Range("E34").Value = '  -0.96%  '
$ws.Range("E35").Value = '  -2.86%  '
$ws.Range("E36").Value = '  -0.63%  '
$ws.Range("E37").Value = '  +4.24%  '
$ws.Range("E38").Value = '  -1.20%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0162'
$ws.Range("E39").Value = '  +0.56%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.94'
$ws.Range("E40").Value = '  +0.71%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.517'
$ws.Range("E41").Value = '  +1.08%  '
$ws.Range("E42").Value = '  -0.51%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.772'
$ws.Range("E43").Value = '  +0.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0463'
$ws.Range("E44").Value = '  -0.43%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '63.34'
$ws.Range("E45").Value = '  +2.86%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.26'
$ws.Range("E46").Value = '  -2.25%  '
$ws.Range("D47").Value = '1.697.07'
$ws.Range("E47").Value = '  +0.29%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.828'
$ws.Range("E48").Value = '  -8.88%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '84.52'
$ws.Range("E49").Value = '  -0.77%  '
$ws.Range("E50").Value = '  +0.04%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '41.55'
$ws.Range("E51").Value = '  -1.42%  '
